$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ---
$ws.Range("A5").Value = "Admin"
$ws.Range("B5").Value = "Admin"
$ws.Range("C5").Value = "IND_DAU_51"
$ws.Range("D5").Value = "IDM+18"
$ws.Range("E5").Value = "10.75.58.51"
$ws.Range("F5").Value = 409026540
$ws.Range("G5").Value = "100"
$ws.Range("H5").Value = "530"
$ws.Range("I5").Value = "60"
$ws.Range("J5").Value = "1"
$ws.Range("K5").Value = "500"

# --- Row 6 ---
$ws.Range("A6").Value = "Admin"
$ws.Range("B6").Value = "Admin"
$ws.Range("C6").Value = "IND_DAU_51"
$ws.Range("D6").Value = "IDM+18"
$ws.Range("E6").Value = "10.75.58.51"
$ws.Range("F6").Value = 409026540
$ws.Range("G6").Value = "100"
$ws.Range("H6").Value = "530"
$ws.Range("I6").Value = "60"
$ws.Range("J6").Value = "1"
$ws.Range("K6").Value = "30000"

# Copy formatting (style incl. quote-prefix) from row 2 onto the new rows,
# matching the existing data rows' look and feel.
$ws.Range("A2:K2").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)
$ws.Range("A2:K2").Copy()
$ws.Range("A6:K6").PasteSpecial(-4122)

# Re-assert the text values (prefixed so Excel keeps them as text / shared
# strings rather than reinterpreting them as numbers after the format paste).
# F5/F6 are left untouched here - they already hold the correct numeric
# value/style from before the format paste.
$ws.Range("G5").Value = "'100"
$ws.Range("H5").Value = "'530"
$ws.Range("I5").Value = "'60"
$ws.Range("J5").Value = "'1"
$ws.Range("K5").Value = "'500"

$ws.Range("G6").Value = "'100"
$ws.Range("H6").Value = "'530"
$ws.Range("I6").Value = "'60"
$ws.Range("J6").Value = "'1"
$ws.Range("K6").Value = "'30000"

$ws.Range("K6").Select()
